$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hourly cryptos-list refresh (coinranking.com) -- GitHub Actions run, Mon May 22 07:21:48 UTC 2023.
# Price (D) / 1h-volume (E) cells are plain text in this sheet (not numbers), and rows 12/13,
# 20/21, 43/44 and 50/51 swap rank order along with their data. Column D values that parse as a
# clean number are briefly forced to Text format so Excel does not silently convert them to
# numerics (which would also rewrite e.g. "5.330" -> "5.33"); the style is restored to Normal
# immediately afterwards so no stray formatting is left behind.

$ws.Range("D2").Value = "27.009.91"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "1.822.03"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.74%  "
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4643"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3665"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07242"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8593"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.04%  "
$ws.Range("E11").Value = "  -3.10%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07623"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.85%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.858.15"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.330"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.500"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.008"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008639"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.007"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("B20").Value = "WrappedBTC"
$ws.Range("C20").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D20").Value = "26.918.68"
$ws.Range("E20").Value = "  -2.32%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.154"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.22%  "
$ws.Range("D24").Value = "2.011.16"
$ws.Range("E24").Value = "  -3.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.841"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.047"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.105"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08836"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.952"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.429"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.129"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7184"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.075"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05258"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01926"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.403"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.928"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.145"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5174"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.63%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1628"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.00%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8595"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -14.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.163"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4807"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.22%  "
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "102.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06251"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.619"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.19%  "
